$d = $word.ActiveDocument

# 1) Freelance web-dev bullet: drop "creation tool" wording, end the sentence with a period.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Designed, implemented, and maintained functional websites for clients using the WordPress creation tool",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Designed, implemented, and maintained functional websites for clients using WordPress.",
    2)
if (-not $found1) { throw "Could not find bullet 1 (WordPress)" }
Write-Output "replace1: $found1"

# 2) Drafted Analytic reports bullet: "to the" -> "for the"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Drafted Analytic reports and visualizations to the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Drafted Analytic reports and visualizations for the ",
    2)
if (-not $found2) { throw "Could not find bullet 2 (Drafted Analytic reports)" }
Write-Output "replace2: $found2"

# 3) "Worked on creating an Email web API" -> "Worked on ASP.NET and SQL to develop the review system."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "Worked on creating an Email web API",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Worked on ASP.NET and SQL to develop the review system.",
    2)
if (-not $found3) { throw "Could not find bullet 3 (Worked on creating an Email web API)" }
Write-Output "replace3: $found3"

# 4) "Applied MVC architecture to develop the Review system" -> "Hands-on experience with creating an Email Web API"
$rng4 = $d.Content
$found4 = $rng4.Find.Execute(
    "Applied MVC architecture to develop the Review system",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Hands-on experience with creating an Email Web API",
    2)
if (-not $found4) { throw "Could not find bullet 4 (Applied MVC architecture)" }
Write-Output "replace4: $found4"

Write-Output "All replacements applied successfully."
